$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 daily spot price update (automatic electricity price refresh)
$ws.Range("A2").Value = 45995
$ws.Range("B2").Value = 61.08
$ws.Range("C2").Value = 55.95
$ws.Range("D2").Value = 47.8
$ws.Range("E2").Value = 43.04
$ws.Range("F2").Value = 35.7
$ws.Range("G2").Value = 46.5
$ws.Range("H2").Value = 61.56
$ws.Range("I2").Value = 77.22
$ws.Range("J2").Value = 84.45999999999999
$ws.Range("K2").Value = 77.93000000000001
$ws.Range("L2").Value = 65.98
$ws.Range("M2").Value = 57.99
$ws.Range("N2").Value = 51.6
$ws.Range("O2").Value = 47.94
$ws.Range("P2").Value = 47.33
$ws.Range("Q2").Value = 53.3
$ws.Range("R2").Value = 65.61
$ws.Range("S2").Value = 79.11
$ws.Range("T2").Value = 85.51000000000001
$ws.Range("U2").Value = 88.95999999999999
$ws.Range("V2").Value = 88.84
$ws.Range("W2").Value = 85.44
$ws.Range("X2").Value = 80.23
$ws.Range("Y2").Value = 65.34999999999999
$ws.Range("Z2").Value = 64.77
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 79.95999999999999
$ws.Range("AD2").Value = 87.23999999999999
$ws.Range("AF2").Value = 87.14
$ws.Range("AG2").Value = "0h-15h"
